# Uloha5/charakteristika.xlsx - "Tabulky a grafy done"
#
# The measurement table header in C3 used to read
#   "Napeti na vystupu zesilovace" + newline + "rezistoru"
# The trailing word "rezistoru" was a leftover/typo and is removed,
# leaving just "Napeti na vystupu zesilovace" followed by a newline.
#
# The two rows (700 kHz / 800 kHz) whose phase measurement could not be
# taken were stored with a nonsensical placeholder number
# (9.8999999999999899E+37). They are replaced with a simple "-" marker,
# matching how the rest of the sheet denotes "no data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "no data" placeholder numbers with a literal "-" first, so
# that once the old C3 header string becomes unused it is dropped from
# the shared string table and the newly introduced strings end up
# ordered "-" before the updated header text (matching how Excel itself
# lays out freshly appended shared strings).
$ws.Range("E35").Value = "-"
$ws.Range("E36").Value = "-"

# Fix up the C3 header text (drop the stray "rezistoru" line).
$ws.Range("C3").Value = "Napětí na výstupu zesilovače" + [char]10
